$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "test - Copy")
$ws.Name = "test1"

# New warning / label cells for the duplicate run_parameters table check
$ws.Range("F1").Value = "Duplicate run_parameters table"

$ws.Range("F2").Value = "Output Path"
$ws.Range("G2").Value = "Version"

# F3 is entered with a leading quote (text/quote-prefix), mirrors source data
$ws.Range("F3").Value = "'output"
$ws.Range("G3").Value = "vTest"

# Add a third table (duplicate run_parameters table) over F2:G3
$null = $ws.ListObjects.Add(1, $ws.Range("F2:G3"), 0, 1)
$tbl3 = $ws.ListObjects.Item(1)
$tbl3.Name = "run_parameters30"
$tbl3.TableStyle = "TableStyleMedium5"

# Match the resulting selection state
$null = $ws.Range("I15").Select()
